$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the "last modified" tracking columns added to the table.
$ws.Range("F1").Value = "Last modified by"
$ws.Range("G1").Value = "Last modified date & time"

# F1 carries a small 9pt font (matches the rest of the "Last modified" metadata style).
$f1Font = $ws.Range("F1").Font
$f1Font.Color = 0
$f1Font.Name = "__Inter_aaf875"
$f1Font.Size = 9

# Widen the two new columns.
$ws.Range("F1").ColumnWidth = 19.83
$ws.Range("G1").ColumnWidth = 22.3

# Scroll / selection state, as left by the author after adding the columns.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("H11").Select()
